$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Legislature" row entirely (base-parameterization row cleanup);
# subsequent rows shift up by one.
$ws.Rows.Item(12).Delete()

# Re-parameterize "Water Rights Division (SWRCB)" row (now row 16):
# was C=1, D=1, E=1, F=1  ->  now C=-0.5, D=-0.5, E and F cleared
$ws.Range("C16").Value = -0.5
$ws.Range("D16").Value = -0.5
$ws.Range("E16").ClearContents()
$ws.Range("F16").ClearContents()

# Update view state to match the saved selection.
$ws.Range("F16").Select()
